$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J (match formatting of existing header H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-19
$values = @{
    2  = @(4, 4)
    3  = @(6, 6)
    4  = @(9, 9)
    5  = @(5, 6)
    6  = @(5, 5)
    7  = @(8, 8)
    8  = @(6, 7)
    9  = @(9, 9)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(7, 8)
    15 = @(8, 9)
    16 = @(5, 5)
    17 = @(5, 6)
    18 = @(7, 7)
    19 = @(5, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
